$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $new
    }
}

# Title
Replace-Text "Unveiling Quantum Entanglement: A Mysterious Dance of Particles" "The Journey Through the Realm of Biology: Unveiling the Secrets of Life"

# Author name
Replace-Text " Emily Carter" " Sophia Kennedy"

# Email address (collapse 5 runs -> 3 runs' worth of text)
Replace-Text "edu" "org"
Replace-Text "emily" "sophiakennedy@educonnect"
Replace-Text "carter@quantamverse." ""

# Intro paragraph
Replace-Text "In the realm of quantum mechanics, a realm where the laws of classical physics falter, lies a perplexing phenomenon known as quantum entanglement" "Delving into the realm of biology is akin to embarking on an extraordinary voyage through the intricate web of life, where mysteries unfold and wonders abound"
Replace-Text " This enigmatic connection between particles transcends the constraints of time and space, allowing them to share information instantaneously, regardless of the distance separating them" " Biology, the study of living organisms, unveils the captivating secrets that orchestrate the symphony of life on our planet"
Replace-Text " In this essay, we will delve into the world of quantum entanglement, unraveling its complexities and delving into its profound implications for our understanding of the universe" " It is a boundless realm where curiosity and discovery converge, beckoning us to explore the enigmatic tapestry of nature's designs"

# Second paragraph (after first <w:br/><w:br/>)
Replace-Text "The strange and elusive nature of quantum entanglement has captivated the minds of scientists and philosophers alike since its discovery in the mid-20th century" "Biology unveils the intricate machinery that governs the very essence of life"
Replace-Text " It defies our intuition and challenges our notions of locality and causality" " From the microscopic symphony of cells to the awe-inspiring majesty of organisms, biology uncovers the profound interconnectedness that binds all living entities"
Replace-Text " As we explore this fascinating phenomenon, we will examine the experiments that have confirmed its existence, the theories that attempt to explain it, and the potential applications that it may hold for the future of computing, cryptography, and information transfer" " It delves into the hidden realms of genetics, where the blueprint of life is inscribed, and explores the marvels of evolution, where organisms adapt, thrive, and transform"

# Third paragraph (after second <w:br/><w:br/>)
Replace-Text "Quantum entanglement has the potential to fundamentally alter our understanding of the universe" "The study of biology transcends mere knowledge acquisition; it fosters an appreciation for the diversity and unity of life"
Replace-Text " If particles can communicate instantaneously over vast distances, it raises profound questions about the nature of reality and the role of locality in the laws of physics" " Through biology, we gain insights into our own existence, unraveling the complexities of human anatomy and physiology, and fostering a profound respect for the delicate balance that sustains life on Earth"
Replace-Text " Furthermore, the ability to manipulate and harness quantum entanglement could open up new avenues for technology, revolutionizing communication, computation, and cryptography. Unveiling the mysteries of quantum entanglement is a scientific endeavor of immense importance, with the potential to reshape our understanding of the universe and transform the way we live" " It is a subject that ignites our imagination, kindles our curiosity, and inspires us to ponder the greatest mysteries of existence"

# Summary heading paragraph
Replace-Text "Quantum entanglement, a mysterious phenomenon in the realm of quantum mechanics, defies our classical understanding of locality and causality" "Biology unveils the intricate tapestry of life, delving into the secrets of living organisms, from the microscopic to the macroscopic"
Replace-Text " Two entangled particles, regardless of their distance apart, share information instantaneously" " It encompasses the study of genetics, evolution, and the interconnectedness of life, fostering an appreciation for diversity and unity"
Replace-Text ". Scientists have conducted experiments confirming the existence of this phenomenon, and theories have emerged to explain its enigmatic characteristics. The potential applications of quantum entanglement are vast, including secure communication, enhanced computing, and the development of new materials" ""

$cr = [char]13
Replace-Text " Unveiling the mysteries of quantum entanglement is a captivating and transformative scientific pursuit, holding the key to unlocking the secrets of the universe and revolutionizing technology." (" Biology ignites curiosity, inspires exploration, and cultivates a profound understanding of our own existence and the intricate web of life on Earth." + $cr)
